# repull data, push all data, mean calculation
# Update column F (dSF) values on Sheet1 to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F3").Value = 4
$ws.Range("F4").Value = -3
$ws.Range("F5").Value = -5
$ws.Range("F6").Value = -7
$ws.Range("F11").Value = -1
$ws.Range("F12").Value = -3
$ws.Range("F14").Value = -1
$ws.Range("F16").Value = -4
$ws.Range("F20").Value = 6
$ws.Range("F24").Value = -2
